# Add a new "JWT" worksheet after the existing "FRONTEND" sheet, documenting
# the JWT / Spring-Security work items, then make it the active sheet.

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the last existing sheet (FRONTEND) so it becomes
# the 4th / last tab.
$sheetCount  = $wb.Worksheets.Count
$lastSheet   = $wb.Worksheets.Item($sheetCount)
$ws          = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name     = "JWT"

# Fill in the cell contents. The write order below (not strict row order)
# matches how the entries were actually authored.
$ws.Range("A3").Value  = "pom.xml"
$ws.Range("A2").Value  = "MumsApplication.java"
$ws.Range("A4").Value  = "application.properties"
$ws.Range("A5").Value  = "config package"
$ws.Range("A6").Value  = "security package"
$ws.Range("A7").Value  = "LoginRepo"
$ws.Range("A8").Value  = "Login "
$ws.Range("A9").Value  = "AppConstants"
$ws.Range("A16").Value = "SalesmanServiceImpl"
$ws.Range("A10").Value = "AuthController"
$ws.Range("A11").Value = "JwtAuthResponse- payload"

# Widen column A so the longer labels are readable.
$ws.Columns.Item(1).ColumnWidth = 35.25

# Make the new sheet the active one, with A11 selected (matches the saved
# view state of the workbook).
[void]$ws.Activate()
[void]$ws.Range("A11").Select()
